$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell without Excel re-interpreting
# numeric-looking strings (e.g. "125.10") as numbers, and without leaving any
# permanent NumberFormat/style change on the cell (matches source cells which
# carry no explicit style).
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue $ws.Range("D2") "54.016.91"
Set-TextValue $ws.Range("E2") "  -8.24%  "
Set-TextValue $ws.Range("D3") "2.857.53"
Set-TextValue $ws.Range("E3") "  -11.04%  "
Set-TextValue $ws.Range("E4") "  -0.07%  "
Set-TextValue $ws.Range("D5") "471.41"
Set-TextValue $ws.Range("E5") "  -11.86%  "
Set-TextValue $ws.Range("D6") "125.10"
Set-TextValue $ws.Range("E6") "  -7.46%  "
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("D8") "2.850.96"
Set-TextValue $ws.Range("E8") "  -11.35%  "
Set-TextValue $ws.Range("D9") "0.400"
Set-TextValue $ws.Range("E9") "  -12.52%  "
Set-TextValue $ws.Range("D10") "6.61"
Set-TextValue $ws.Range("E10") "  -12.63%  "
Set-TextValue $ws.Range("D11") "0.0955"
Set-TextValue $ws.Range("E11") "  -16.24%  "
Set-TextValue $ws.Range("D12") "0.328"
Set-TextValue $ws.Range("E12") "  -16.48%  "
Set-TextValue $ws.Range("E13") "  -4.66%  "
Set-TextValue $ws.Range("D14") "3.343.47"
Set-TextValue $ws.Range("E14") "  -10.99%  "
Set-TextValue $ws.Range("D15") "22.88"
Set-TextValue $ws.Range("E15") "  -11.58%  "
Set-TextValue $ws.Range("D16") "54.015.62"
Set-TextValue $ws.Range("E16") "  -8.30%  "
Set-TextValue $ws.Range("D17") "2.861.82"
Set-TextValue $ws.Range("E17") "  -10.81%  "
Set-TextValue $ws.Range("D18") "0.0000133"
Set-TextValue $ws.Range("E18") "  -15.75%  "
Set-TextValue $ws.Range("D19") "5.29"
Set-TextValue $ws.Range("E19") "  -10.69%  "
Set-TextValue $ws.Range("D20") "11.43"
Set-TextValue $ws.Range("E20") "  -13.98%  "
Set-TextValue $ws.Range("D21") "7.01"
Set-TextValue $ws.Range("E21") "  -14.55%  "
Set-TextValue $ws.Range("D22") "291.35"
Set-TextValue $ws.Range("E22") "  -19.22%  "
Set-TextValue $ws.Range("E23") "  -0.17%  "
Set-TextValue $ws.Range("D24") "0.441"
Set-TextValue $ws.Range("E24") "  -14.97%  "
Set-TextValue $ws.Range("D25") "58.42"
Set-TextValue $ws.Range("E25") "  -16.70%  "
Set-TextValue $ws.Range("E26") "  +0.23%  "
Set-TextValue $ws.Range("E27") "  -11.31%  "
Set-TextValue $ws.Range("E28") "  -0.17%  "
Set-TextValue $ws.Range("D29") "0.0₃0803"
Set-TextValue $ws.Range("E29") "  -16.73%  "
Set-TextValue $ws.Range("D30") "6.22"
Set-TextValue $ws.Range("E30") "  -12.33%  "
Set-TextValue $ws.Range("D31") "1.11"
Set-TextValue $ws.Range("E31") "  -7.92%  "
Set-TextValue $ws.Range("D32") "6.11"
Set-TextValue $ws.Range("E32") "  -13.40%  "
Set-TextValue $ws.Range("D33") "18.87"
Set-TextValue $ws.Range("E33") "  -13.14%  "
Set-TextValue $ws.Range("D34") "1.60"
Set-TextValue $ws.Range("E34") "  -16.73%  "
Set-TextValue $ws.Range("D35") "4.18"
Set-TextValue $ws.Range("E35") "  -14.70%  "
Set-TextValue $ws.Range("D36") "134.65"
Set-TextValue $ws.Range("E36") "  -16.75%  "
Set-TextValue $ws.Range("D37") "5.38"
Set-TextValue $ws.Range("E37") "  -15.46%  "
Set-TextValue $ws.Range("E38") "  -16.37%  "
Set-TextValue $ws.Range("D39") "22.71"
Set-TextValue $ws.Range("E39") "  -12.59%  "
Set-TextValue $ws.Range("D40") "2.877.78"
Set-TextValue $ws.Range("E40") "  -11.09%  "
Set-TextValue $ws.Range("E41") "  -13.71%  "
Set-TextValue $ws.Range("E42") "  -0.17%  "
Set-TextValue $ws.Range("D43") "35.03"
Set-TextValue $ws.Range("E43") "  -14.49%  "
Set-TextValue $ws.Range("D44") "0.953"
Set-TextValue $ws.Range("E44") "  -13.10%  "
Set-TextValue $ws.Range("D45") "0.595"
Set-TextValue $ws.Range("E45") "  -16.71%  "
Set-TextValue $ws.Range("D46") "3.37"
Set-TextValue $ws.Range("E46") "  -16.24%  "
Set-TextValue $ws.Range("E47") "  -13.34%  "
Set-TextValue $ws.Range("D48") "2.032.59"
Set-TextValue $ws.Range("E48") "  -11.60%  "
Set-TextValue $ws.Range("D49") "5.26"
Set-TextValue $ws.Range("E49") "  -16.12%  "
Set-TextValue $ws.Range("D50") "17.63"
Set-TextValue $ws.Range("E50") "  -15.11%  "
Set-TextValue $ws.Range("D51") "0.0209"
Set-TextValue $ws.Range("E51") "  -12.68%  "
